$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the employee name above the existing clock entries.
$ws.Range("A1").Value = "Freddy Velez"

# Fix a typo in the clock-in id (8FD8ADBD20 -> 4FD889D140).
$ws.Range("A2").Value = "4FD889D140"

# A3 ("IN -> 2017/01/31 18:57") is unchanged.

# Append the corresponding clock-out entry.
$ws.Range("A4").Value = "OUT -> 2017/02/14 17:22"
